$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the updated "Price" (column D) values are plain decimal numbers
# (e.g. "0.9980", "15.90") that must stay stored as text, exactly like every other
# cell in this column, instead of being auto-coerced into numbers by Excel (which
# would silently drop significant trailing zeros, e.g. "0.9980" -> 0.998). Mark just
# those specific cells as Text before writing so only their value -- not their style --
# changes, matching every other (already-text) cell in the sheet.
$textValueCells = @("D4", "D5", "D6", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textValueCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.722.40'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '1.730.98'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("D4").Value = '0.9983'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '242.33'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '0.9984'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").Value = '0.4929'
$ws.Range("E7").Value = '  +1.26%  '
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").Value = '0.06217'
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").Value = '1.727.79'
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("D11").Value = '15.90'
$ws.Range("E11").Value = '  +3.16%  '
$ws.Range("D12").Value = '0.06997'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = '0.6115'
$ws.Range("E13").Value = '  +2.25%  '
$ws.Range("D14").Value = '4.505'
$ws.Range("E14").Value = '  -0.70%  '
$ws.Range("D15").Value = '77.25'
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").Value = '0.9980'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").Value = '26.526.50'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = '0.9987'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("D19").Value = '0.000007222'
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("D21").Value = '1.948.16'
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("D22").Value = '4.492'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").Value = '8.580'
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '5.107'
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("D25").Value = '138.51'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").Value = '15.37'
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("D27").Value = '1.773'
$ws.Range("E27").Value = '  +3.19%  '
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("D29").Value = '106.51'
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D31").Value = '0.07991'
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").Value = '3.677'
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").Value = '0.9973'
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.608'
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.003'
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.6235'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '0.9414'
$ws.Range("E38").Value = '  +3.72%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.042'
$ws.Range("E39").Value = '  +2.53%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.422'
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '0.9980'
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.01511'
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.595'
$ws.Range("E43").Value = '  +3.76%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '99.35'
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3860'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '6.938'
$ws.Range("E46").Value = '  +3.69%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1160'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05385'
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '7.857'
$ws.Range("E49").Value = '  +2.06%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '30.24'
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '51.79'
$ws.Range("E51").Value = '  +1.51%  '
